# Test-file audit cleanup: remove the stray "Sheet" placeholder row that had
# been left in the optimization_parameters sheet (A16:C16 = "Sheet", 3, 4).
# Deleting the whole row shifts "simulation_timepoints" (previously row 17)
# up into row 16, and drops the now-unused "Sheet" shared string / numFmt
# style automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Activate() | Out-Null
$ws.Rows.Item(16).Select() | Out-Null
$ws.Rows.Item(16).Delete() | Out-Null
